$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# Rename "Done" header to "Done in %"
$ws.Range("D1").Value = "Done in %"

# Update the selection/active cell to D2
$ws.Range("D2").Select()
